$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").NumberFormat = "mm-dd-yy"
$ws.Range("H1").Value = (Get-Date -Year 2015 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

$ws.Range("H2").Value = " "
$ws.Range("H3").Value = " "
$ws.Range("H4").Value = " "
$ws.Range("H5").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("H7").Value = " "

$ws.Range("H7").Select()
